$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, pushing existing data down to row 3
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the new verse data
$ws.Range("A2").Value = "40860_lxx"
$ws.Range("B2").Value = "40860_mt"
$ws.Range("C2").Value = "ἐὰν δὲ ὁ κύριος ᾖ μετ’ αὐτοῦ, οὐκ ἀποτίσει · ἐὰν δὲ μισθωτὸς ᾖ, ἔσται αὐτῷ ἀντὶ τοῦ μισθοῦ αὐτοῦ."
# D2 stays blank (Greek Clause is empty for this verse)
$ws.Range("E2").Value = "(22, 14)"
# F2 stays blank (no Greek preposition for this verse)
$ws.Range("G2").Value = ">M B<LJW <MW L> JCLM >M FKJR HW> B> B FKRW"
$ws.Range("H2").Value = "B FKRW"
$ws.Range("I2").Value = "prep"
$ws.Range("J2").Value = "B"
$ws.Range("K2").Value = "inanim"
